$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the current row 13 ("*After reviewing..." notice row)
$ws.Rows.Item(13).Insert()

# B13: label "Total"
$ws.Range("B13").Value = "សរុប"
$ws.Range("B13").HorizontalAlignment = -4152

# C13:F13 sum formulas
$ws.Range("C13").Formula = "=SUM(C11:C12)"
$ws.Range("D13").Formula = "=SUM(D11:D12)"
$ws.Range("E13").Formula = "=SUM(E11:E12)"
$ws.Range("F13").Formula = "=SUM(F11:F12)"

Write-Output "done"
